# Populate the "2023" worksheet (sheet5) with weekend demand data
# and mark it as the active/selected sheet (replacing "2022" as active).

$wb = $excel.ActiveWorkbook
$ws2022 = $wb.Worksheets.Item("2022")
$ws2023 = $wb.Worksheets.Item("2023")

# Header label in B1 ("Demand")
$ws2023.Range("B1").Value = "Demand"

# Weekend date-range labels (column A) and demand values (column B)
$data = @(
  "07Jan-08Jan|110768.79309203901",
  "14Jan-15Jan|110180.35862964259",
  "21Jan-22Jan|101037.35365170331",
  "28Jan-29Jan|102569.8714283099",
  "04Feb-05Feb|103652.5880219853",
  "11Feb-12Feb|105489.6587806874",
  "18Feb-19Feb|106550.2612862548",
  "25Feb-26Feb|106770.0676989444",
  "04Mar-05Mar|99242.402498299518",
  "11Mar-12Mar|102989.79566850699",
  "18Mar-19Mar|105266.8702056771",
  "25Mar-26Mar|104277.88173043181",
  "01Apr-02Apr|104507.391589005",
  "08Apr-09Apr|103803.36416431741",
  "15Apr-16Apr|107893.21880058511",
  "22Apr-23Apr|107728.7567571343",
  "29Apr-30Apr|105305.68813733679",
  "06May-07May|106528.2878560813",
  "13May-14May|109206.86660674569",
  "20May-21May|105244.1554141605",
  "27May-28May|106359.6940331059",
  "03Jun-04Jun|103005.6602175969",
  "10Jun-11Jun|107532.9148656321",
  "17Jun-18Jun|105510.38955361499",
  "24Jun-25Jun|105215.4267140382",
  "01Jul-02Jul|104423.1932100478",
  "08Jul-09Jul|105101.22683635241",
  "15Jul-16Jul|103729.8596960111",
  "22Jul-23Jul|103859.0053991366",
  "29Jul-30Jul|104595.9628707378",
  "05Aug-06Aug|106574.5571561523",
  "12Aug-13Aug|105928.4054998764",
  "19Aug-20Aug|103480.7593371425",
  "26Aug-27Aug|105079.972774037",
  "02Sep-03Sep|105374.3734722899",
  "09Sep-10Sep|102849.0898070931",
  "16Sep-17Sep|106943.42025951399",
  "23Sep-24Sep|103575.89899272579",
  "30Sep-01Oct|103298.5258705218",
  "07Oct-08Oct|110182.02171206981",
  "14Oct-15Oct|109447.7106363498",
  "21Oct-22Oct|104720.05875072371",
  "28Oct-29Oct|105939.3168219752",
  "04Nov-05Nov|107924.8897817336",
  "11Nov-12Nov|105906.7896748292",
  "18Nov-19Nov|106847.73781658171",
  "25Nov-26Nov|104565.6777242027",
  "02Dec-03Dec|105411.2916825182",
  "09Dec-10Dec|108067.3287371935",
  "16Dec-17Dec|106068.00830065169",
  "23Dec-24Dec|104430.11714667649",
  "30Dec-31Dec|103861.7363291918"
)

$row = 2
foreach ($item in $data) {
    $parts = $item.Split("|")
    $label = $parts[0]
    $val = [double]$parts[1]
    $ws2023.Cells.Item($row, 1).Value = $label
    $ws2023.Cells.Item($row, 2).Value = $val
    $row = $row + 1
}
$lastRow = $row - 1

# Match the formatting (bold, centered, bordered) already used for the
# equivalent label cells on the other yearly sheets.
$ws2022.Range("B1").Copy()
$ws2023.Range("B1").PasteSpecial(-4122)

$ws2022.Range("A2").Copy()
$ws2023.Range("A2:A" + $lastRow).PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Make "2023" the active sheet (was "2022"), with E4 selected.
$ws2023.Activate()
$ws2023.Range("E4").Select()
